# Update "Förändrad" (column C) date value from 45184 to 45186 for every
# data row, and add the record's "Beteckning" (column A) as the friendly
# name (2nd argument) of every HYPERLINK() formula in columns S-Y that
# does not already have one.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count + $usedRange.Row - 1

$oldDate = 45184
$newDate = 45186

# Columns that may contain HYPERLINK(...) formulas needing a friendly name.
$hyperlinkCols = @(19, 20, 21, 22, 23, 24, 25)  # S, T, U, V, W, X, Y

for ($r = 2; $r -le $lastRow; $r++) {

    # --- Column C: bump the "Förändrad" date value ---
    $cCell = $ws.Cells.Item($r, 3)
    if ($cCell.Value2 -eq $oldDate) {
        $cCell.Value2 = $newDate
    }

    # --- Columns S..Y: append friendly name to HYPERLINK formulas ---
    $beteckning = $ws.Cells.Item($r, 1).Value2

    if ($beteckning) {
        foreach ($col in $hyperlinkCols) {
            $cell = $ws.Cells.Item($r, $col)
            if ($cell.HasFormula) {
                $formula = $cell.Formula
                if ($formula.Contains("HYPERLINK(") -and -not $formula.Contains(",")) {
                    $trimmed = $formula.Substring(0, $formula.Length - 1)
                    $cell.Formula = $trimmed + ', "' + $beteckning + '")'
                }
            }
        }
    }
}
